# "Cohérence des dimensions du bras et beaux vérins"
# Update the dimensioning parameters of the excavator arm ("Feuil1") so
# the bucket/arm/cylinder geometry stays consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Bucket tooth / bucket dimensions -------------------------------------
$ws.Range("B1").Value = 70                 # Longueur_dent_godet
$ws.Range("B2").Value = 50                 # Largeur_dent_godet (was =B6/6, now a plain value)
$ws.Range("B3").Value = 10                 # Epaisseur_dent_godet
$ws.Range("B4").Formula = "=B1*2/5"        # Depassement_dent_godet (now derived from B1)
$ws.Range("B5").Value = 600                # Longueur_godet
$ws.Range("B6").Value = 450                # Largeur_godet
$ws.Range("B7").Value = 200                # Rayon_courbure_godet
$ws.Range("B8").Value = 42                 # Angle_dattaque_godet
$ws.Range("B9").Value = 20                 # Epaisseur_godet
$ws.Range("B12").Value = 150               # Longeur_liaison_godet
$ws.Range("B13").Value = 100               # Epaisseur_liaison_godet

# --- Forearm / forearm cylinder -------------------------------------------
$ws.Range("B30").Value = 1800              # Longueur_avant_bras
$ws.Range("B31").Formula = "=0.9*B30"      # Longueur_verin_avant_bras (ratio 0.94 -> 0.9)

# --- Arm cylinders ("beaux vérins") ----------------------------------------
$ws.Range("B53").Value = 2000              # Longueur_verin_bras
$ws.Range("B54").Value = 1500              # Longueur_verin_chassis
$ws.Range("B55").Value = 140               # Angle_bras
$ws.Range("B56").Value = 150               # Largeur_verin_bras
$ws.Range("B60").Formula = "=MAX(B53/15,B59*2)"  # Rayon_ext_avant_bras_bras
$ws.Range("B67").Formula = "=B54*2/3"      # Longueur_verin_bras_chassis (now derived)
$ws.Range("B68").Formula = "=B54*2/3"      # Longueur_tige_verin_bras_chassis (now derived)

# --- Restore the cursor position left by the author on save ----------------
$ws.Range("B69").Select()
